# Adds the 16-Mar / 17-Mar daily-track entries (rows 36-42) to the MAR-22
# sheet, mirroring the existing formatting conventions used elsewhere in
# the workbook (copy formats from analogous existing cells so styles.xml
# reuses identical cellXfs entries rather than minting near-duplicates).

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("MAR-22")
$ws2 = $wb.Worksheets.Item("FEB-22")

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Row 36 - "16" record header (uses the "no outer border" family of
# styles, like row 30 of FEB-22).
# ---------------------------------------------------------------------
$ws3.Cells.Item(36, 1).Value2 = 16

$ws2.Cells.Item(30, 2).Copy() | Out-Null
$ws3.Cells.Item(36, 2).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(36, 2).Value2 = 44642

$ws2.Cells.Item(30, 6).Copy() | Out-Null
$ws3.Cells.Item(36, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(36, 3).Value2 = "RPA GSS"

$ws2.Cells.Item(30, 6).Copy() | Out-Null
$ws3.Cells.Item(36, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(36, 4).WrapText = $true
$ws3.Cells.Item(36, 4).Value2 = "1. Implementation of Public holidays at Activity Customer_token system task is work in progress"

$ws2.Cells.Item(30, 5).Copy() | Out-Null
$ws3.Cells.Item(36, 5).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(36, 5).Value2 = 0.5

$ws2.Cells.Item(30, 6).Copy() | Out-Null
$ws3.Cells.Item(36, 6).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(36, 6).Value2 = "WIP"

# ---------------------------------------------------------------------
# Row 37 - continuation line (bordered-blank family, as row 34).
# ---------------------------------------------------------------------
$ws3.Cells.Item(34, 1).Copy() | Out-Null
$ws3.Cells.Item(37, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 2).Copy() | Out-Null
$ws3.Cells.Item(37, 2).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 3).Copy() | Out-Null
$ws3.Cells.Item(37, 3).PasteSpecial($xlPasteFormats) | Out-Null

$ws2.Cells.Item(29, 4).Copy() | Out-Null
$ws3.Cells.Item(37, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(37, 4).Value2 = "2. Correction Received for the duplicate record at database (as task running twice) is work in progress"

$ws3.Cells.Item(34, 5).Copy() | Out-Null
$ws3.Cells.Item(37, 5).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(37, 5).Value2 = 0.1

$ws2.Cells.Item(29, 6).Copy() | Out-Null
$ws3.Cells.Item(37, 6).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(37, 6).Value2 = "WIP"

# ---------------------------------------------------------------------
# Row 38 - continuation line, "Completed".
# ---------------------------------------------------------------------
$ws3.Cells.Item(34, 1).Copy() | Out-Null
$ws3.Cells.Item(38, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 2).Copy() | Out-Null
$ws3.Cells.Item(38, 2).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 3).Copy() | Out-Null
$ws3.Cells.Item(38, 3).PasteSpecial($xlPasteFormats) | Out-Null

$ws2.Cells.Item(29, 4).Copy() | Out-Null
$ws3.Cells.Item(38, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(38, 4).Value2 = "3. For the RLOGIC, General ledger has been triggered for all three centers of the Feb22"

$ws3.Cells.Item(34, 5).Copy() | Out-Null
$ws3.Cells.Item(38, 5).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(38, 5).Value2 = 1

$ws2.Cells.Item(29, 6).Copy() | Out-Null
$ws3.Cells.Item(38, 6).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(38, 6).Value2 = "Completed"

# ---------------------------------------------------------------------
# Row 39 - "17" record header (standard bordered family, as row 31/33).
# ---------------------------------------------------------------------
$ws3.Cells.Item(31, 1).Copy() | Out-Null
$ws3.Cells.Item(39, 1).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(39, 1).Value2 = 17

$ws3.Cells.Item(31, 2).Copy() | Out-Null
$ws3.Cells.Item(39, 2).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(39, 2).Value2 = 44643

$ws3.Cells.Item(31, 3).Copy() | Out-Null
$ws3.Cells.Item(39, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(39, 3).Value2 = "RPA GSS"

$ws3.Cells.Item(31, 4).Copy() | Out-Null
$ws3.Cells.Item(39, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(39, 4).Value2 = "1. Correction received from the GRS-Details task due to multiple time database has been updated for single running, " + [char]10 + "Now it has been fixed, tested and it is running smoothly and a few logs are implemented"

$ws3.Cells.Item(31, 5).Copy() | Out-Null
$ws3.Cells.Item(39, 5).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(39, 5).Value2 = 1

$ws3.Cells.Item(31, 6).Copy() | Out-Null
$ws3.Cells.Item(39, 6).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(39, 6).Value2 = "Completed"

$ws3.Rows.Item(39).RowHeight = 28.8

# ---------------------------------------------------------------------
# Row 40 - continuation line (bordered-blank family, as row 34), tall row.
# ---------------------------------------------------------------------
$ws3.Cells.Item(34, 1).Copy() | Out-Null
$ws3.Cells.Item(40, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 2).Copy() | Out-Null
$ws3.Cells.Item(40, 2).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 3).Copy() | Out-Null
$ws3.Cells.Item(40, 3).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(32, 4).Copy() | Out-Null
$ws3.Cells.Item(40, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(40, 4).Value2 = "2. Public holidays has been implemented at Customer Visit token task, it has been tested and it is running smoothly, whereas" + [char]10 + "a few captcha correction works is work in progress"

$ws3.Cells.Item(34, 5).Copy() | Out-Null
$ws3.Cells.Item(40, 5).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(40, 5).Value2 = 0.3

$ws3.Cells.Item(34, 6).Copy() | Out-Null
$ws3.Cells.Item(40, 6).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(40, 6).Value2 = "WIP"

$ws3.Rows.Item(40).RowHeight = 43.2

# ---------------------------------------------------------------------
# Row 41 - continuation line (bordered-blank family, as row 34).
# ---------------------------------------------------------------------
$ws3.Cells.Item(34, 1).Copy() | Out-Null
$ws3.Cells.Item(41, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 2).Copy() | Out-Null
$ws3.Cells.Item(41, 2).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 3).Copy() | Out-Null
$ws3.Cells.Item(41, 3).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(34, 4).Copy() | Out-Null
$ws3.Cells.Item(41, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(41, 4).Value2 = "3. Activity Daily has been received to enhance the task running logic and it is work in progress for multiple centers"

$ws3.Cells.Item(34, 5).Copy() | Out-Null
$ws3.Cells.Item(41, 5).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(41, 5).Value2 = 0.5

$ws3.Cells.Item(34, 6).Copy() | Out-Null
$ws3.Cells.Item(41, 6).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(41, 6).Value2 = "WIP"

# ---------------------------------------------------------------------
# Row 42 - trailing "Comments"-style note row (as row 30).
# ---------------------------------------------------------------------
$ws3.Cells.Item(30, 1).Copy() | Out-Null
$ws3.Cells.Item(42, 1).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(30, 2).Copy() | Out-Null
$ws3.Cells.Item(42, 2).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(30, 3).Copy() | Out-Null
$ws3.Cells.Item(42, 3).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(42, 3).Value2 = "RPA RLOGIC"

$ws3.Cells.Item(30, 4).Copy() | Out-Null
$ws3.Cells.Item(42, 4).PasteSpecial($xlPasteFormats) | Out-Null
$ws3.Cells.Item(42, 4).Value2 = "4. Rlogic,  email issue exists and yet to be fixed."

$ws3.Cells.Item(30, 5).Copy() | Out-Null
$ws3.Cells.Item(42, 5).PasteSpecial($xlPasteFormats) | Out-Null

$ws3.Cells.Item(30, 6).Copy() | Out-Null
$ws3.Cells.Item(42, 6).PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------
# Reposition the view the way MS Excel left it after the edit.
# ---------------------------------------------------------------------
$ws3.Application.Goto($ws3.Range("A31"), $true)
$ws3.Range("D44").Select() | Out-Null
